# =====================================================================
# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. directly
#    before the existing "2022-Q2" sheet), carrying the same layout/
#    styling as the other quarterly sheets, and fill it with the fund
#    holdings for 2022-Q3.
# 2. Insert a new row into the "总计" summary sheet for 2022-Q3,
#    pushing the existing quarters down by one row.
#
# NOTE: all structural row/sheet operations (Insert / full-row
# PasteSpecial) are performed FIRST. Only afterwards is a scratch cell
# set up for literal-text staging, because a whole-row Insert() would
# otherwise propagate the scratch cell's Text number-format across the
# newly inserted row (and inflate the sheet's used range) if the
# scratch cell already existed.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1a. Create the new "2022-Q3" worksheet by duplicating the current
#     "2022-Q2" sheet (position 2) so that all sheet-level formatting
#     (sheetPr, column widths, cell styles, page margins, etc.) carries
#     over exactly, then rename it; it naturally lands right before
#     the original "2022-Q2" sheet.
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item(2)
$srcSheet.Copy($srcSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The duplicated template only has 7 data rows (rows 2-8) but 2022-Q3
# needs 8 data rows (2-9), so append 2 more rows, copying the row-7
# formatting down so the new rows match the existing style.
$newSheet.Range("A7:H7").Copy()
$newSheet.Range("A8:H8").PasteSpecial(-4122)
$newSheet.Range("A7:H7").Copy()
$newSheet.Range("A9:H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 1b. Update the "总计" summary sheet: insert a new row for 2022-Q3
#     directly below the header row, shifting the existing quarters
#     down by one row.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Range("A3:D3").Copy()
$ws1.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Helper: write $text into $range as a literal text value (no
#    locale number-sniffing, no left-over number-format/quote-prefix
#    residue on the destination cell) by staging it through a scratch
#    cell pre-formatted as Text, then copying the VALUE ONLY into
#    place. Set up now that no more whole-row Insert() calls remain.
# ---------------------------------------------------------------------
$scratch = $ws1.Range("ZZ1")
$scratch.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

function Set-Row($sheet, $r, $a, $b, $c, $d, $e, $f, $g, $h) {
    $sheet.Range("A" + $r).Value = $a
    Set-TextValue $sheet.Range("B" + $r) $b
    Set-TextValue $sheet.Range("C" + $r) $c
    Set-TextValue $sheet.Range("D" + $r) $d
    Set-TextValue $sheet.Range("E" + $r) $e
    Set-TextValue $sheet.Range("F" + $r) $f
    Set-TextValue $sheet.Range("G" + $r) $g
    $sheet.Range("H" + $r).Value = $h
}

# ---------------------------------------------------------------------
# 3. Fill the new "2022-Q3" sheet with the fund holdings data.
# ---------------------------------------------------------------------
Set-Row $newSheet 2 0 "005258" "景顺长城量化平衡灵活配置混合" "1.17" "90.11" "2.61" "0.0305" 5
Set-Row $newSheet 3 1 "015061" "中信建投沪深300指数增强A" "1.47" "91.17" "1.72" "0.0253" 4
Set-Row $newSheet 4 2 "001244" "华泰柏瑞量化智慧灵活配置混合A" "2.90" "92.68" "0.76" "0.0220" 4
Set-Row $newSheet 5 3 "015062" "中信建投沪深300指数增强C" "1.14" "91.17" "1.72" "0.0196" 4
Set-Row $newSheet 6 4 "006063" "景顺长城MSCI中国A股国际通指数增强" "0.58" "93.59" "2.07" "0.0120" 10
Set-Row $newSheet 7 5 "014861" "申万菱信双禧混合A" "1.51" "30.74" "0.58" "0.0088" 1
Set-Row $newSheet 8 6 "006104" "华泰柏瑞量化智慧灵活配置混合C" "0.38" "92.68" "0.76" "0.0029" 4
Set-Row $newSheet 9 7 "014862" "申万菱信双禧混合C" "0.04" "30.74" "0.58" "0.0002" 1

# ---------------------------------------------------------------------
# 4. Fill the new 2022-Q3 row on the "总计" sheet, then renumber the
#    running index in column A for all the rows that shifted down.
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0.12

for ($r = 3; $r -le 8; $r++) {
    $ws1.Cells.Item($r, 1).Value = ($r - 2)
}

# ---------------------------------------------------------------------
# 5. Clean up the scratch cell used for text staging.
# ---------------------------------------------------------------------
$scratch.Clear()
